$d = $word.ActiveDocument

# Insert a "." between the closing "<<hearingType>>" merge field and the
# following "<<cs_{hearingType==...}>>" conditional merge field so the
# sentence "The hearing will be <<hearingType>>." reads correctly.
$d.Content.Find.Execute(">><<cs_", $true, $false, $false, $false, $false, `
                         $true, 1, $false, ">>.<<cs_", 2)
